$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format first so that numeric-looking
# strings (e.g. "29.209.12", "1.000", "5.223") are not auto-converted
# into numbers/dates by Excel when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.209.12'
$ws.Range("E2").Value = '  -0.98%  '

# Row 3
$ws.Range("D3").Value = '1.866.59'
$ws.Range("E3").Value = '  -0.60%  '

# Row 4
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '0.7096'
$ws.Range("E5").Value = '  -0.93%  '

# Row 6
$ws.Range("D6").Value = '241.78'
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '0.3111'
$ws.Range("E8").Value = '  +0.10%  '

# Row 9
$ws.Range("D9").Value = '0.07651'
$ws.Range("E9").Value = '  -4.17%  '

# Row 10
$ws.Range("D10").Value = '24.68'
$ws.Range("E10").Value = '  -2.73%  '

# Row 11
$ws.Range("D11").Value = '0.08368'
$ws.Range("E11").Value = '  +0.98%  '

# Row 12
$ws.Range("D12").Value = '1.866.28'
$ws.Range("E12").Value = '  -0.21%  '

# Row 13
$ws.Range("D13").Value = '5.223'
$ws.Range("E13").Value = '  -1.16%  '

# Row 14
$ws.Range("D14").Value = '0.7105'
$ws.Range("E14").Value = '  -2.75%  '

# Row 15
$ws.Range("D15").Value = '91.30'
$ws.Range("E15").Value = '  +0.11%  '

# Row 16
$ws.Range("D16").Value = '29.218.96'
$ws.Range("E16").Value = '  -0.98%  '

# Row 17
$ws.Range("D17").Value = '5.955'
$ws.Range("E17").Value = '  +0.31%  '

# Row 18
$ws.Range("D18").Value = '243.58'
$ws.Range("E18").Value = '  -0.87%  '

# Row 19
$ws.Range("D19").Value = '0.000007827'
$ws.Range("E19").Value = '  -0.70%  '

# Row 20
$ws.Range("D20").Value = '2.115.68'
$ws.Range("E20").Value = '  +0.07%  '

# Row 21
$ws.Range("E21").Value = '  -2.08%  '

# Row 22
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").Value = '7.860'
$ws.Range("E23").Value = '  -1.50%  '

# Row 24
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  +0.07%  '

# Row 25
$ws.Range("D25").Value = '0.1630'
$ws.Range("E25").Value = '  +1.25%  '

# Row 26
$ws.Range("D26").Value = '163.35'
$ws.Range("E26").Value = '  -0.13%  '

# Row 27
$ws.Range("E27").Value = '  -1.27%  '

# Row 28
$ws.Range("E28").Value = '  +0.95%  '

# Row 29
$ws.Range("E29").Value = '  +0.81%  '

# Row 30
$ws.Range("E30").Value = '  -3.56%  '

# Row 31
$ws.Range("D31").Value = '4.397'
$ws.Range("E31").Value = '  +0.10%  '

# Row 32
$ws.Range("D32").Value = '4.247'
$ws.Range("E32").Value = '  +3.06%  '

# Row 33
$ws.Range("D33").Value = '0.05146'
$ws.Range("E33").Value = '  -2.32%  '

# Row 34
$ws.Range("D34").Value = '0.7954'
$ws.Range("E34").Value = '  +9.21%  '

# Row 35
$ws.Range("D35").Value = '1.912'
$ws.Range("E35").Value = '  -2.59%  '

# Row 36
$ws.Range("E36").Value = '  -2.80%  '

# Row 37
$ws.Range("D37").Value = '2.686'
$ws.Range("E37").Value = '  +0.26%  '

# Row 38
$ws.Range("D38").Value = '0.01856'

# Row 39
$ws.Range("E39").Value = '  -0.21%  '

# Row 40
$ws.Range("D40").Value = '1.157.42'
$ws.Range("E40").Value = '  -5.50%  '

# Row 41
$ws.Range("D41").Value = '6.328'
$ws.Range("E41").Value = '  +3.29%  '

# Row 42
$ws.Range("D42").Value = '0.8970'
$ws.Range("E42").Value = '  -1.69%  '

# Row 43
$ws.Range("E43").Value = '  -0.93%  '

# Row 44
$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.09%  '

# Row 45
$ws.Range("D45").Value = '103.36'
$ws.Range("E45").Value = '  +1.06%  '

# Row 46
$ws.Range("D46").Value = '2.011.57'
$ws.Range("E46").Value = '  -0.13%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.5142'
$ws.Range("E47").Value = '  -2.77%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.777'
$ws.Range("E48").Value = '  -1.49%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  -0.87%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.333'
$ws.Range("E50").Value = '  -0.16%  '

# Row 51
$ws.Range("D51").Value = '0.4290'

# Restore the default (Normal) style on the Price column so the cell
# formatting matches the original workbook (no explicit style).
$ws.Range("D2:D51").Style = "Normal"
